$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 769.63635
$ws.Range("I38").Value = 274
$ws.Range("J38").Value = 3000
$ws.Range("K38").Value = 822
$ws.Range("L38").Value = 9000
$ws.Range("M38").Value = -450
$ws.Range("N38").Value = -9744

$ws.Range("H51").Value = 3029.9333
$ws.Range("I51").Value = 2494.4443
$ws.Range("J51").Value = 3833.1667
$ws.Range("K51").Value = 2494.4443
$ws.Range("L51").Value = 3833.1667
$ws.Range("M51").Value = -2010.4443
$ws.Range("N51").Value = -4801.1667

$ws.Range("H92").Value = 731.3333
$ws.Range("I92").Value = 586.34485
$ws.Range("J92").Value = 1782.5
$ws.Range("K92").Value = 586.34485
$ws.Range("L92").Value = 1782.5
$ws.Range("M92").Value = 661.65515
$ws.Range("N92").Value = -4278.5

$ws.Range("H111").Value = 775.8
$ws.Range("I111").Value = 493
$ws.Range("J111").Value = 1200
$ws.Range("K111").Value = 1479
$ws.Range("L111").Value = 3600
$ws.Range("M111").Value = 1588
$ws.Range("N111").Value = -9734

$ws.Range("H125").Value = 21220
$ws.Range("I125").Value = 34100
$ws.Range("K125").Value = 306900
$ws.Range("M125").Value = -304440

$ws.Range("H135").Value = 3800.3809
$ws.Range("I135").Value = 927.9231
$ws.Range("J135").Value = 8468.125
$ws.Range("K135").Value = 8351.3079
$ws.Range("L135").Value = 76213.125
$ws.Range("M135").Value = -5816.3079
$ws.Range("N135").Value = -81283.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 924.5454999999999
$ws.Range("I110").Value = 968
$ws.Range("K110").Value = 968
$ws.Range("M110").Value = 1077

$ws.Range("H119").Value = 47341.43
$ws.Range("J119").Value = 47341.43
$ws.Range("L119").Value = 47341.43
$ws.Range("N119").Value = -57017.43

$ws.Range("H120").Value = 32000
$ws.Range("J120").Value = 32000
$ws.Range("L120").Value = 32000
$ws.Range("N120").Value = -41676

$ws.Range("H132").Value = 154175.83
$ws.Range("I132").Value = 3717.2
$ws.Range("K132").Value = 11151.6
$ws.Range("M132").Value = -8621.599999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1611.3572
$ws.Range("I107").Value = 1471.5834
$ws.Range("K107").Value = 1471.5834
$ws.Range("M107").Value = 448.4166

$ws.Range("H134").Value = 21880.826
$ws.Range("I134").Value = 25630.092
$ws.Range("J134").Value = 1259.875
$ws.Range("K134").Value = 76890.276
$ws.Range("L134").Value = 3779.625
$ws.Range("M134").Value = -74355.276
$ws.Range("N134").Value = -8849.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 20000
$ws.Range("J29").Value = 20000
$ws.Range("L29").Value = 20000
$ws.Range("N29").Value = -20586

$ws.Range("H31").Value = 17059.148
$ws.Range("I31").Value = 25255.604
$ws.Range("J31").Value = 1544.4286
$ws.Range("K31").Value = 25255.604
$ws.Range("L31").Value = 1544.4286
$ws.Range("M31").Value = -24960.604
$ws.Range("N31").Value = -2134.4286

$ws.Range("H34").Value = 17059.148
$ws.Range("I34").Value = 25255.604
$ws.Range("J34").Value = 1544.4286
$ws.Range("K34").Value = 25255.604
$ws.Range("L34").Value = 1544.4286
$ws.Range("M34").Value = -25053.604
$ws.Range("N34").Value = -1948.4286

$ws.Range("H99").Value = 1235.238
$ws.Range("I99").Value = 1149.875
$ws.Range("J99").Value = 1508.4
$ws.Range("K99").Value = 1149.875
$ws.Range("L99").Value = 1508.4
$ws.Range("M99").Value = 348.125
$ws.Range("N99").Value = -4504.4

$ws.Range("H107").Value = 1178.7826
$ws.Range("I107").Value = 1288.3125
$ws.Range("J107").Value = 928.4286
$ws.Range("K107").Value = 1288.3125
$ws.Range("L107").Value = 928.4286
$ws.Range("M107").Value = 631.6875
$ws.Range("N107").Value = -4768.4286

$ws.Range("H126").Value = 1235.238
$ws.Range("I126").Value = 1149.875
$ws.Range("J126").Value = 1508.4
$ws.Range("K126").Value = 3449.625
$ws.Range("L126").Value = 4525.200000000001
$ws.Range("M126").Value = -979.625
$ws.Range("N126").Value = -9465.200000000001

$ws.Range("H132").Value = 1855.9565
$ws.Range("I132").Value = 1287.5883
$ws.Range("J132").Value = 3466.3333
$ws.Range("K132").Value = 3862.7649
$ws.Range("L132").Value = 10398.9999
$ws.Range("M132").Value = -1332.7649
$ws.Range("N132").Value = -15458.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 773.9474
$ws.Range("I92").Value = 856
$ws.Range("J92").Value = 700.1
$ws.Range("K92").Value = 2568
$ws.Range("L92").Value = 2100.3
$ws.Range("M92").Value = -1320
$ws.Range("N92").Value = -4596.3

$ws.Range("H131").Value = 1673333.9
$ws.Range("I131").Value = 1457.3334
$ws.Range("J131").Value = 2080006.5
$ws.Range("K131").Value = 4372.0002
$ws.Range("L131").Value = 6240019.5
$ws.Range("M131").Value = 667.9997999999996
$ws.Range("N131").Value = -6250099.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 7680
$ws.Range("J23").Value = 36400
$ws.Range("L23").Value = 36400
$ws.Range("N23").Value = -36846

$ws.Range("H24").Value = 16538.5
$ws.Range("I24").Value = 600
$ws.Range("J24").Value = 21851.334
$ws.Range("K24").Value = 600
$ws.Range("L24").Value = 21851.334
$ws.Range("M24").Value = -427
$ws.Range("N24").Value = -22197.334

$ws.Range("H80").Value = 2927.8
$ws.Range("I80").Value = 2879
$ws.Range("J80").Value = 2940
$ws.Range("K80").Value = 2879
$ws.Range("L80").Value = 2940
$ws.Range("M80").Value = -1881
$ws.Range("N80").Value = -4936

$ws.Range("H83").Value = 2927.8
$ws.Range("I83").Value = 2879
$ws.Range("J83").Value = 2940
$ws.Range("K83").Value = 14395
$ws.Range("L83").Value = 14700
$ws.Range("M83").Value = -9403
$ws.Range("N83").Value = -24684

$ws.Range("H107").Value = 965.6667
$ws.Range("I107").Value = 697.1429000000001
$ws.Range("J107").Value = 1341.6
$ws.Range("K107").Value = 697.1429000000001
$ws.Range("L107").Value = 1341.6
$ws.Range("M107").Value = 1222.8571
$ws.Range("N107").Value = -5181.6

$ws.Range("H126").Value = 1246.1666
$ws.Range("I126").Value = 1218.25
$ws.Range("J126").Value = 1302
$ws.Range("K126").Value = 3654.75
$ws.Range("L126").Value = 3906
$ws.Range("M126").Value = -1184.75
$ws.Range("N126").Value = -8846

$ws.Range("H132").Value = 10784.333
$ws.Range("I132").Value = 8257.875
$ws.Range("J132").Value = 15837.25
$ws.Range("K132").Value = 24773.625
$ws.Range("L132").Value = 47511.75
$ws.Range("M132").Value = -22243.625
$ws.Range("N132").Value = -52571.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 50387.273
$ws.Range("J108").Value = 50387.273
$ws.Range("L108").Value = 50387.273
$ws.Range("N108").Value = -58067.273

$ws.Range("H119").Value = 17816.666
$ws.Range("J119").Value = 17816.666
$ws.Range("L119").Value = 17816.666
$ws.Range("N119").Value = -27492.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 765000
$ws.Range("J15").Value = 765000
$ws.Range("L15").Value = 765000
$ws.Range("N15").Value = -765576

$ws.Range("H96").Value = 2001333.4
$ws.Range("I96").Value = 1000000
$ws.Range("J96").Value = 2502000
$ws.Range("K96").Value = 1000000
$ws.Range("L96").Value = 2502000
$ws.Range("M96").Value = -998627
$ws.Range("N96").Value = -2504746

$ws.Range("H107").Value = 285
$ws.Range("I107").Value = 117.9
$ws.Range("J107").Value = 619.2
$ws.Range("K107").Value = 353.7
$ws.Range("L107").Value = 1857.6
$ws.Range("M107").Value = 1566.3
$ws.Range("N107").Value = -5697.6

$ws.Range("H119").Value = 44398.855
$ws.Range("J119").Value = 44398.855
$ws.Range("L119").Value = 44398.855
$ws.Range("N119").Value = -54074.855

$ws.Range("H122").Value = 772915.5600000001
$ws.Range("I122").Value = 2003400.6
$ws.Range("J122").Value = 3862.375
$ws.Range("K122").Value = 6010201.800000001
$ws.Range("L122").Value = 11587.125
$ws.Range("M122").Value = -6007751.800000001
$ws.Range("N122").Value = -16487.125
